$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N ("Late"), shifting Late/heading/Outstanding
# one column to the right (N->O, O->P, P->Q), matching the new "Variable Instalments"
# column added to the repayment schedule.
$ws.Columns("N:N").Insert(-4161) | Out-Null

# Match the width of the newly inserted column to the "In Advance" column (M) width,
# which results in a stored width of 11 (custom, not best-fit).
$ws.Columns("N:N").ColumnWidth = $ws.Columns("M:M").ColumnWidth

# Make "Repayment schedule" the active sheet/tab, and select cell Q7 on it (this also
# clears the previously active selection/tab on the "NewLoanInput" sheet).
$ws.Activate() | Out-Null
$ws.Range("Q7").Select() | Out-Null
